# Apply the two content edits described by the diff:
#  1. Notes Master "datetimeFigureOut" date field: 10/1/2020 -> 10/20/2024
#  2. Typo fix on Slide 5 (Body Forces), Content Placeholder 2:
#     "A gravitational body force can be through of as the " ->
#     "A gravitational body force can be thought of as the "

$p = $ppt.ActivePresentation

# --- 1. Update the date field shown on the Notes Master ---
$notesMaster = $p.NotesMaster
$dateHf = $notesMaster.HeadersFooters.DateAndTime
$dateHf.Text = "10/20/2024"

# --- 2. Fix "through of" -> "thought of" typo on slide 5 ---
$oldText = "A gravitational body force can be through of as the "
$newText = "A gravitational body force can be thought of as the "

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count
            for ($pi = 1; $pi -le $paraCount; $pi++) {
                $para = $tr.Paragraphs($pi)
                $runCount = $para.Runs().Count
                for ($ri = 1; $ri -le $runCount; $ri++) {
                    $run = $para.Runs($ri)
                    if ($run.Text -eq $oldText) {
                        $run.Text = $newText
                    }
                }
            }
        }
    }
}
